$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "tets2"

$ws.Activate()
$ws.Range("B2").Select()
